$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("H11").Value = 122.8
$ws.Range("I11").Value = 122.8
$ws.Range("K11").Value = 122.8
$ws.Range("M11").Value = 17.2
$ws.Range("H125").Value = 7997.5
$ws.Range("J125").Value = 7997
$ws.Range("L125").Value = 71973
$ws.Range("N125").Value = -76893
$ws.Range("H137").Value = 3099.375
$ws.Range("I137").Value = 1900
$ws.Range("J137").Value = 3499.1667
$ws.Range("K137").Value = 5700
$ws.Range("L137").Value = 10497.5001
$ws.Range("M137").Value = -3150
$ws.Range("N137").Value = -15597.5001
$ws.Range("H138").Value = 6735.385
$ws.Range("J138").Value = 7584.6875
$ws.Range("L138").Value = 22754.0625
$ws.Range("N138").Value = -33034.0625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 999.5
$ws.Range("I2").Value = 999
$ws.Range("K2").Value = 999
$ws.Range("M2").Value = -886
$ws.Range("H32").Value = 2443.6978
$ws.Range("I32").Value = 1729.0256
$ws.Range("J32").Value = 9411.75
$ws.Range("K32").Value = 1729.0256
$ws.Range("L32").Value = 9411.75
$ws.Range("M32").Value = -1442.0256
$ws.Range("N32").Value = -9985.75
$ws.Range("H61").Value = 3428.3
$ws.Range("I61").Value = 3285.75
$ws.Range("K61").Value = 3285.75
$ws.Range("M61").Value = -3073.75
$ws.Range("H74").Value = 1188.4546
$ws.Range("I74").Value = 730.3333
$ws.Range("K74").Value = 730.3333
$ws.Range("M74").Value = 143.6667
$ws.Range("H77").Value = 1188.4546
$ws.Range("I77").Value = 730.3333
$ws.Range("K77").Value = 3651.6665
$ws.Range("M77").Value = 716.3334999999997
$ws.Range("H110").Value = 3452.1428
$ws.Range("I110").Value = 3528.9167
$ws.Range("J110").Value = 2991.5
$ws.Range("K110").Value = 3528.9167
$ws.Range("L110").Value = 2991.5
$ws.Range("M110").Value = -1483.9167
$ws.Range("N110").Value = -7081.5
$ws.Range("H116").Value = 999.5
$ws.Range("I116").Value = 999
$ws.Range("K116").Value = 999
$ws.Range("M116").Value = 1295
$ws.Range("H136").Value = 3428.3
$ws.Range("I136").Value = 3285.75
$ws.Range("K136").Value = 9857.25
$ws.Range("M136").Value = -7307.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 999.5
$ws.Range("I3").Value = 999
$ws.Range("K3").Value = 999
$ws.Range("M3").Value = -885
$ws.Range("H134").Value = 3858.1765
$ws.Range("I134").Value = 3790
$ws.Range("J134").Value = 4949
$ws.Range("K134").Value = 11370
$ws.Range("L134").Value = 14847
$ws.Range("M134").Value = -8835
$ws.Range("N134").Value = -19917

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("H3").Value = 17997
$ws.Range("I3").Value = 17997
$ws.Range("K3").Value = 17997
$ws.Range("M3").Value = -17884
$ws.Range("H16").Value = 8778
$ws.Range("I16").Value = 8778
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 8778
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -8491
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 5710.619
$ws.Range("I31").Value = 3732.625
$ws.Range("J31").Value = 6927.846
$ws.Range("K31").Value = 3732.625
$ws.Range("L31").Value = 6927.846
$ws.Range("M31").Value = -3437.625
$ws.Range("N31").Value = -7517.846
$ws.Range("H34").Value = 5710.619
$ws.Range("I34").Value = 3732.625
$ws.Range("J34").Value = 6927.846
$ws.Range("K34").Value = 3732.625
$ws.Range("L34").Value = 6927.846
$ws.Range("M34").Value = -3530.625
$ws.Range("N34").Value = -7331.846
$ws.Range("H113").Value = 8778
$ws.Range("I113").Value = 8778
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 8778
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -6608
$ws.Range("N113").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item(5)
$ws.Range("H18").Value = 154.5
$ws.Range("I18").Value = 154.5
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 463.5
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -294.5
$ws.Range("N18").ClearContents()
$ws.Range("H68").Value = 1436.6
$ws.Range("J68").Value = 1495.75
$ws.Range("L68").Value = 4487.25
$ws.Range("N68").Value = -6109.25
$ws.Range("H71").Value = 1436.6
$ws.Range("J71").Value = 1495.75
$ws.Range("L71").Value = 13461.75
$ws.Range("N71").Value = -21573.75
$ws.Range("H107").Value = 3327.5
$ws.Range("I107").Value = 4664.4
$ws.Range("K107").Value = 13993.2
$ws.Range("M107").Value = -12073.2
$ws.Range("H123").Value = 5837.2
$ws.Range("I123").Value = 2744
$ws.Range("K123").Value = 8232
$ws.Range("M123").Value = -5782

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 42.4
$ws.Range("I2").Value = 28
$ws.Range("K2").Value = 28
$ws.Range("M2").Value = 85
$ws.Range("H3").Value = 4671
$ws.Range("I3").Value = 2177.5
$ws.Range("K3").Value = 2177.5
$ws.Range("M3").Value = -2061.5
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H11").Value = 7602000.5
$ws.Range("I11").Value = 8502501
$ws.Range("K11").Value = 8502501
$ws.Range("M11").Value = -8502362
$ws.Range("H13").Value = 205
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H80").Value = 4393.722
$ws.Range("I80").Value = 3684
$ws.Range("J80").Value = 5103.4443
$ws.Range("K80").Value = 3684
$ws.Range("L80").Value = 5103.4443
$ws.Range("M80").Value = -2686
$ws.Range("N80").Value = -7099.4443
$ws.Range("H83").Value = 4393.722
$ws.Range("I83").Value = 3684
$ws.Range("J83").Value = 5103.4443
$ws.Range("K83").Value = 18420
$ws.Range("L83").Value = 25517.2215
$ws.Range("M83").Value = -13428
$ws.Range("N83").Value = -35501.2215
$ws.Range("H113").Value = 1394.6666
$ws.Range("J113").Value = 1392
$ws.Range("L113").Value = 1392
$ws.Range("N113").Value = -5732

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 4658.8
$ws.Range("J40").Value = 4550
$ws.Range("L40").Value = 4550
$ws.Range("N40").Value = -4822
$ws.Range("H82").Value = 1104.1428
$ws.Range("I82").Value = 956.6667
$ws.Range("K82").Value = 956.6667
$ws.Range("M82").Value = -595.6667
$ws.Range("H85").Value = 1104.1428
$ws.Range("I85").Value = 956.6667
$ws.Range("K85").Value = 956.6667
$ws.Range("M85").Value = 291.3333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item(8)
$ws.Range("H4").Value = 15345.818
$ws.Range("J4").Value = 9271.429
$ws.Range("L4").Value = 9271.429
$ws.Range("N4").Value = -9497.429
$ws.Range("H7").Value = 1000
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -887
$ws.Range("N7").ClearContents()
$ws.Range("H8").Value = 1000
$ws.Range("J8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("N8").Value = -1280
$ws.Range("H10").Value = 13332.333
$ws.Range("I10").Value = 9998.5
$ws.Range("K10").Value = 9998.5
$ws.Range("M10").Value = -9829.5
$ws.Range("H13").Value = 5333.3335
$ws.Range("I13").Value = 3000
$ws.Range("K13").Value = 3000
$ws.Range("M13").Value = -2860
$ws.Range("H14").Value = 466.66666
$ws.Range("I14").Value = 466.66666
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 466.66666
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -298.66666
$ws.Range("N14").ClearContents()
$ws.Range("H17").Value = 10500
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 20000
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = -828
$ws.Range("N17").Value = -20344
$ws.Range("H81").Value = 3750.4666
$ws.Range("I81").Value = 4011.5715
$ws.Range("J81").Value = 95
$ws.Range("K81").Value = 8023.143
$ws.Range("L81").Value = 190
$ws.Range("M81").Value = -6962.143
$ws.Range("N81").Value = -2312
$ws.Range("H84").Value = 3750.4666
$ws.Range("I84").Value = 4011.5715
$ws.Range("J84").Value = 95
$ws.Range("K84").Value = 40115.715
$ws.Range("L84").Value = 950
$ws.Range("M84").Value = -34811.715
$ws.Range("N84").Value = -11558
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 786.1667
$ws.Range("I113").Value = 844.5
$ws.Range("K113").Value = 2533.5
$ws.Range("M113").Value = -363.5
$ws.Range("H122").Value = 500
$ws.Range("I122").Value = 500
$ws.Range("K122").Value = 1500
$ws.Range("M122").Value = 950
$ws.Range("H132").Value = 2895.3103
$ws.Range("I132").Value = 2568.6086
$ws.Range("K132").Value = 7705.825800000001
$ws.Range("M132").Value = -5175.825800000001
